# Applies the "add vertex 42 imagery to copyright sheet" edit.
#
# Summary of the change:
#  - On the copyright sheet ("©", sheet1), insert a new blank row above the
#    current row 3 ("By Vertex42.com"). This shifts all following content
#    rows down by one. The new row gets a custom height of 37 (presumably
#    to host the new Vertex42 logo/imagery).
#  - The copyright sheet becomes the active/selected sheet (tabSelected),
#    with the active cell now at B13 (the "Please review..." license row).
#  - The "template" sheet is no longer the selected/active tab, and its
#    stored active-cell selection moves to D20.
#  - The workbook-level active tab moves from the template sheet (index 1)
#    back to the copyright sheet (index 0).

$wb = $excel.ActiveWorkbook

$copyrightSheet = $wb.Worksheets.Item("©")
$templateSheet = $wb.Worksheets.Item("template")

# Insert a new blank row above row 3 on the copyright sheet; this pushes the
# existing rows 3-16 down to rows 4-17.
$copyrightSheet.Rows.Item(3).Insert()

# The newly inserted row 3 is blank and gets a custom row height of 37.
$copyrightSheet.Rows.Item(3).RowHeight = 37

# Make the copyright sheet the active sheet/tab, with the selection on B13.
$copyrightSheet.Activate()
$copyrightSheet.Range("B13").Select()

# The template sheet keeps a stored selection, now at D20, but is no longer
# the active tab.
$templateSheet.Range("D20").Select()

# Re-activate the copyright sheet so it ends up as the workbook's active tab.
$copyrightSheet.Activate()
